$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: PathFile
$ws.Range("E1").Value = "PathFile"
$ws.Range("E2").Value = "int32"
for ($r = 3; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = 9002
}

# Apply vertical-center alignment to the new column's cells
$ws.Range("E1:E12").VerticalAlignment = -4108

# Update the active selection as in the authored workbook
$ws.Range("G6").Select() | Out-Null
